$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.270.57'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '2.249.95'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.45'
$ws.Range("E5").Value = '  -1.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.93'
$ws.Range("E6").Value = '  -1.22%  '

$ws.Range("E7").Value = '  +1.07%  '

$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.29'
$ws.Range("E10").Value = '  -1.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0817'
$ws.Range("E11").Value = '  -0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.30'
$ws.Range("E12").Value = '  -0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'

$ws.Range("D14").Value = '2.595.34'
$ws.Range("E14").Value = '  +1.08%  '

$ws.Range("D15").Value = '2.329.79'
$ws.Range("E15").Value = '  +4.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.837'
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.66'
$ws.Range("E17").Value = '  -2.44%  '

$ws.Range("D18").Value = '44.137.24'
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("D19").Value = '0.0₃0973'
$ws.Range("E19").Value = '  +1.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.24'
$ws.Range("E20").Value = '  -6.11%  '

$ws.Range("E21").Value = '  +1.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.72'
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.00'
$ws.Range("E23").Value = '  +0.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.04'
$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.26'
$ws.Range("E28").Value = '  +4.74%  '

$ws.Range("E29").Value = '  +1.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.01'
$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.22'
$ws.Range("E31").Value = '  +1.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.31'
$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("E33").Value = '  -2.62%  '

$ws.Range("E34").Value = '  -0.66%  '

$ws.Range("E35").Value = '  +3.57%  '

$ws.Range("E36").Value = '  +3.14%  '

$ws.Range("E37").Value = '  -0.30%  '

$ws.Range("E38").Value = '  -5.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.69'
$ws.Range("E39").Value = '  +4.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.72'
$ws.Range("E40").Value = '  -4.62%  '

$ws.Range("E41").Value = '  -2.99%  '

$ws.Range("E42").Value = '  -1.83%  '

$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("D44").Value = '1.759.78'
$ws.Range("E44").Value = '  +3.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '83.90'
$ws.Range("E45").Value = '  +0.38%  '

$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.75'
$ws.Range("E47").Value = '  -0.70%  '

$ws.Range("E48").Value = '  -2.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.19'
$ws.Range("E49").Value = '  +1.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.08'
$ws.Range("E50").Value = '  -1.42%  '

$ws.Range("E51").Value = '  -3.68%  '
